# Generate Report for Handback
#
# This mirrors the "Handback" report-generation step: for each language
# sheet (zh-cn, de-de) the two file rows move from "In Translation" to
# "Handed back: in sync with en-US", and the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns (I/J/K)
# are populated with the handed-back file info + a handback timestamp.
# The Overview sheet (and the language sheets) also get a few columns
# widened so the new long file names / status text are readable.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/49bdc166b002be578af95397265bd2b4fac0396f/e2e/"
$mdName1 = "091adbed-9c36-4a54-9925-1526bf5c20a5.md"
$mdName2 = "0de5070b-5d4b-4dba-b706-fea4ede3c02f.md"

# --- Overview sheet: widen the per-language status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ===================== zh-cn =====================
$ws = $wb.Worksheets.Item("zh-cn")

# Widen Status (C), Latest Target File (I) and Latest Handback File (J)
$ws.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws.Columns.Item(10).ColumnWidth = 39.166666666666664

# Row 2 -> 091adbed...
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("J2").Value = $ws.Range("G2").Value2
$ws.Range("K2").Value = "2016-08-13 04:28:39"
$ws.Hyperlinks.Add($ws.Range("I2"), ($baseUrl + $mdName1), "", "", $mdName1)

# Row 3 -> 0de5070b...
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("J3").Value = $ws.Range("G3").Value2
$ws.Range("K3").Value = "2016-08-13 04:28:39"
$ws.Hyperlinks.Add($ws.Range("I3"), ($baseUrl + $mdName2), "", "", $mdName2)

# ===================== de-de =====================
$ws = $wb.Worksheets.Item("de-de")

# Widen Status (C), Latest Target File (I) and Latest Handback File (J)
$ws.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws.Columns.Item(10).ColumnWidth = 39.166666666666664

# Row 2 -> 091adbed...
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("J2").Value = $ws.Range("G2").Value2
$ws.Range("K2").Value = "2016-08-13 04:28:48"
$ws.Hyperlinks.Add($ws.Range("I2"), ($baseUrl + $mdName1), "", "", $mdName1)

# Row 3 -> 0de5070b...
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("J3").Value = $ws.Range("G3").Value2
$ws.Range("K3").Value = "2016-08-13 04:28:48"
$ws.Hyperlinks.Add($ws.Range("I3"), ($baseUrl + $mdName2), "", "", $mdName2)
